# Weekly price update: insert two new rows (newest week's "Primera" and
# "Segunda" entries for Cebollín) at row 255, pushing the rest of the
# data block down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 255:256 - this shifts old rows 255-388 down to
# rows 257-390, matching the new dimension A1:R390.
$ws.Range("A255:A256").EntireRow.Insert()

# New row 255: "Primera" entry for the new week (2022-01-11), reusing the
# most recent "Primera" pricing figures.
$ws.Range("A255").Value = 3
$ws.Range("B255").Value = "Femacal de La Calera"
$ws.Range("C255").Value = "Coquimbo"
$ws.Range("D255").Value = "2022-01-11"
$ws.Range("E255").Value = 5
$ws.Range("F255").Value = 100112037
$ws.Range("G255").Value = "Cebollín"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 160
$ws.Range("K255").Value = 3000
$ws.Range("L255").Value = 3000
$ws.Range("M255").Value = 3000
$ws.Range("N255").Value = '$/paquete 36 unidades'
$ws.Range("O255").Value = "Provincia de Quillota"
$ws.Range("P255").Value = 83
$ws.Range("Q255").Value = 36
$ws.Range("R255").Value = "Hortaliza"

# New row 256: "Segunda" entry for the new week (2022-01-11), reusing the
# most recent "Segunda" pricing figures except for the updated volume.
$ws.Range("A256").Value = 3
$ws.Range("B256").Value = "Femacal de La Calera"
$ws.Range("C256").Value = "Coquimbo"
$ws.Range("D256").Value = "2022-01-11"
$ws.Range("E256").Value = 5
$ws.Range("F256").Value = 100112037
$ws.Range("G256").Value = "Cebollín"
$ws.Range("H256").Value = "Sin especificar"
$ws.Range("I256").Value = "Segunda"
$ws.Range("J256").Value = 110
$ws.Range("K256").Value = 2000
$ws.Range("L256").Value = 2000
$ws.Range("M256").Value = 2000
$ws.Range("N256").Value = '$/paquete 36 unidades'
$ws.Range("O256").Value = "Provincia de Quillota"
$ws.Range("P256").Value = 56
$ws.Range("Q256").Value = 36
$ws.Range("R256").Value = "Hortaliza"
